$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing the exact textual representation).
foreach ($addr in @("D5", "D6", "D8", "D9", "D11", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume values scraped on 2024-11-07.
$ws.Range("D2").Value = '75.866.67'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.891.72'
$ws.Range("E3").Value = '  +7.45%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '196.27'
$ws.Range("E5").Value = '  +3.91%  '
$ws.Range("D6").Value = '600.93'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.554'
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").Value = '0.193'
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("D10").Value = '2.891.47'
$ws.Range("E10").Value = '  +7.51%  '
$ws.Range("D11").Value = '0.403'
$ws.Range("E11").Value = '  +11.75%  '
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("E13").Value = '  +4.19%  '
$ws.Range("D14").Value = '3.435.05'
$ws.Range("E14").Value = '  +7.77%  '
$ws.Range("D15").Value = '75.783.51'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").Value = '0.0000191'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '27.51'
$ws.Range("E17").Value = '  +2.98%  '
$ws.Range("D18").Value = '2.898.92'
$ws.Range("E18").Value = '  +8.09%  '
$ws.Range("D19").Value = '8.99'
$ws.Range("E19").Value = '  -4.92%  '
$ws.Range("D20").Value = '12.61'
$ws.Range("E20").Value = '  +4.08%  '
$ws.Range("D21").Value = '383.19'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = '2.31'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").Value = '72.03'
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '3.043.01'
$ws.Range("E26").Value = '  +7.52%  '
$ws.Range("D27").Value = '4.27'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").Value = '9.83'
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("D29").Value = '0.0000108'
$ws.Range("E29").Value = '  +12.19%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = '1.41'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").Value = '511.94'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").Value = '7.83'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  +2.83%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").Value = '165.47'
$ws.Range("E36").Value = '  +1.45%  '
$ws.Range("D37").Value = '20.22'
$ws.Range("E37").Value = '  +4.37%  '
$ws.Range("E38").Value = '  -3.94%  '
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("D40").Value = '184.08'
$ws.Range("E40").Value = '  +7.25%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '0.346'
$ws.Range("E42").Value = '  +3.92%  '
$ws.Range("D43").Value = '5.04'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.68'
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '0.0920'
$ws.Range("E45").Value = '  +7.88%  '
$ws.Range("D46").Value = '1.23'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").Value = '40.45'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '0.586'
$ws.Range("E49").Value = '  +7.96%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '3.78'
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.668'
$ws.Range("E51").Value = '  +12.27%  '
